# Weekly refresh of the "Fruta / hortaliza" sheet: the date (Fecha) and
# the volume/price columns (Volumen, Precio mínimo, Precio máximo,
# Precio promedio ponderado, Precio $/Kg) for data rows 2-49 are
# reshuffled to a newer daily pull. Every other column (IDs, región,
# categoría, variedad, calidad, unidad, origen, clasificación, etc.) is
# identical for all 48 rows, so the only thing that actually needs to
# move is this 6-column block, row by row, according to the mapping
# below (newRow -> oldRow it now carries).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowMap = @{
    2=29;  3=14;  4=16;  5=27;  6=21;  7=30;  8=44;  9=47;  10=33;
    11=22; 12=20; 13=28; 14=4;  15=11; 16=3;  17=36; 18=35; 19=24;
    20=23; 21=39; 22=2;  23=41; 24=13; 25=32; 26=42; 27=9;  28=19;
    29=48; 30=25; 31=38; 32=34; 33=8;  34=5;  35=31; 36=15; 37=46;
    38=40; 39=43; 40=10; 41=12; 42=17; 43=45; 44=7;  45=26; 46=37;
    47=18; 48=49; 49=6
}

# Column indexes for the block that moves with each row.
$colD = 4   # Fecha
$colJ = 10  # Volumen
$colK = 11  # Precio minimo
$colL = 12  # Precio maximo
$colM = 13  # Precio promedio ponderado
$colP = 16  # Precio $/Kg

# Snapshot every data row BEFORE writing anything, since rows both
# read from and write to this same range (the map is not simply a
# shift - several rows point at each other).
$snapshot = @{}
foreach ($r in 2..49) {
    $snapshot[$r] = @{
        D = $ws.Cells.Item($r, $colD).Value2()
        J = $ws.Cells.Item($r, $colJ).Value2()
        K = $ws.Cells.Item($r, $colK).Value2()
        L = $ws.Cells.Item($r, $colL).Value2()
        M = $ws.Cells.Item($r, $colM).Value2()
        P = $ws.Cells.Item($r, $colP).Value2()
    }
}

foreach ($newRow in 2..49) {
    $oldRow = $rowMap[$newRow]
    $src = $snapshot[$oldRow]

    $ws.Cells.Item($newRow, $colD).Value = $src.D
    $ws.Cells.Item($newRow, $colJ).Value = $src.J
    $ws.Cells.Item($newRow, $colK).Value = $src.K
    $ws.Cells.Item($newRow, $colL).Value = $src.L
    $ws.Cells.Item($newRow, $colM).Value = $src.M
    $ws.Cells.Item($newRow, $colP).Value = $src.P
}
